# Updated BGR model - 2025-08-15 20:17
# Update the "lcoe_class" (column P) rankings on the "solar" and "wind"
# resource sheets to reflect the refreshed cost-class ordering.

$wb = $excel.ActiveWorkbook

# --- solar sheet: spv-BGR_16 cost-class block (rows 4-6) ---
$wsSolar = $wb.Worksheets.Item("solar")
$wsSolar.Range("P4").Value = 4
$wsSolar.Range("P5").Value = 2
$wsSolar.Range("P6").Value = 3

# --- wind sheet: several won-BGR cost-class blocks ---
$wsWind = $wb.Worksheets.Item("wind")

# won-BGR_29 block (rows 4-5)
$wsWind.Range("P4").Value = 2
$wsWind.Range("P5").Value = 3

# won-BGR_25 block (rows 15-16)
$wsWind.Range("P15").Value = 2
$wsWind.Range("P16").Value = 3

# won-BGR_24 block (rows 18-20)
$wsWind.Range("P18").Value = 2
$wsWind.Range("P19").Value = 3
$wsWind.Range("P20").Value = 1

# won-BGR_17 block (rows 47-48)
$wsWind.Range("P47").Value = 1
$wsWind.Range("P48").Value = 2
